# Auto-generated edit script: updates cached numeric values in the
# "Sheets/Zalera_Profits.xlsx" workbook (profit-tracking tables) across
# all 8 worksheets, per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2070.4814
$ws.Range("I135").Value = 1272.3889
$ws.Range("J135").Value = 3666.6667
$ws.Range("K135").Value = 11451.5001
$ws.Range("L135").Value = 33000.0003
$ws.Range("M135").Value = -8916.500099999999
$ws.Range("N135").Value = -38070.0003
$ws.Range("H138").Value = 4178.0405
$ws.Range("J138").Value = 4691.6055
$ws.Range("L138").Value = 14074.8165
$ws.Range("N138").Value = -24354.8165

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 537.125
$ws.Range("I4").Value = 560.5714
$ws.Range("J4").Value = 373
$ws.Range("K4").Value = 560.5714
$ws.Range("L4").Value = 373
$ws.Range("M4").Value = -444.5714
$ws.Range("N4").Value = -605
$ws.Range("H74").Value = 716247.1
$ws.Range("I74").Value = 835121.75
$ws.Range("J74").Value = 2999.5
$ws.Range("K74").Value = 835121.75
$ws.Range("L74").Value = 2999.5
$ws.Range("M74").Value = -834247.75
$ws.Range("N74").Value = -4747.5
$ws.Range("H77").Value = 716247.1
$ws.Range("I77").Value = 835121.75
$ws.Range("J77").Value = 2999.5
$ws.Range("K77").Value = 4175608.75
$ws.Range("L77").Value = 14997.5
$ws.Range("M77").Value = -4171240.75
$ws.Range("N77").Value = -23733.5
$ws.Range("H110").Value = 6758670
$ws.Range("I110").Value = 10001261
$ws.Range("J110").Value = 3272
$ws.Range("K110").Value = 10001261
$ws.Range("L110").Value = 3272
$ws.Range("M110").Value = -9999216
$ws.Range("N110").Value = -7362
$ws.Range("H132").Value = 8728.083000000001
$ws.Range("I132").Value = 3629.2222
$ws.Range("J132").Value = 24024.666
$ws.Range("K132").Value = 10887.6666
$ws.Range("L132").Value = 72073.99800000001
$ws.Range("M132").Value = -8357.6666
$ws.Range("N132").Value = -77133.99800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1998.75
$ws.Range("I22").Value = 2170.8572
$ws.Range("K22").Value = 2170.8572
$ws.Range("M22").Value = -1997.8572
$ws.Range("H86").Value = 119484.766
$ws.Range("I86").Value = 1952.5625
$ws.Range("K86").Value = 1952.5625
$ws.Range("M86").Value = -829.5625
$ws.Range("H89").Value = 119484.766
$ws.Range("I89").Value = 1952.5625
$ws.Range("K89").Value = 9762.8125
$ws.Range("M89").Value = -4146.8125
$ws.Range("H99").Value = 2498.5833
$ws.Range("I99").Value = 2202.8948
$ws.Range("K99").Value = 2202.8948
$ws.Range("M99").Value = -704.8948
$ws.Range("H105").Value = 45469136
$ws.Range("I105").Value = 66686212
$ws.Range("K105").Value = 66686212
$ws.Range("M105").Value = -66684465
$ws.Range("H107").Value = 1852
$ws.Range("I107").Value = 1732.5
$ws.Range("K107").Value = 1732.5
$ws.Range("M107").Value = 187.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1901.5714
$ws.Range("I107").Value = 1005.5
$ws.Range("J107").Value = 2260
$ws.Range("K107").Value = 1005.5
$ws.Range("L107").Value = 2260
$ws.Range("M107").Value = 914.5
$ws.Range("N107").Value = -6100
$ws.Range("H122").Value = 935.6
$ws.Range("I122").Value = 984
$ws.Range("K122").Value = 2952
$ws.Range("M122").Value = -502
$ws.Range("H132").Value = 25004.615
$ws.Range("I132").Value = 3381.182
$ws.Range("K132").Value = 10143.546
$ws.Range("M132").Value = -7613.545999999998
$ws.Range("H134").Value = 4130.6763
$ws.Range("I134").Value = 3814.9395
$ws.Range("J134").Value = 14550
$ws.Range("K134").Value = 11444.8185
$ws.Range("L134").Value = 43650
$ws.Range("M134").Value = -8909.818499999999
$ws.Range("N134").Value = -48720

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2120299.8
$ws.Range("I4").Value = 637197.25
$ws.Range("K4").Value = 1911591.75
$ws.Range("M4").Value = -1911479.75
$ws.Range("H11").Value = 125110.94
$ws.Range("I11").Value = 9448
$ws.Range("K11").Value = 28344
$ws.Range("M11").Value = -28204
$ws.Range("H15").Value = 157.66667
$ws.Range("I15").Value = 158.88889
$ws.Range("J15").Value = 155.83333
$ws.Range("K15").Value = 476.66667
$ws.Range("L15").Value = 467.49999
$ws.Range("M15").Value = -336.66667
$ws.Range("N15").Value = -747.49999
$ws.Range("H17").Value = 782.7143
$ws.Range("J17").Value = 807
$ws.Range("L17").Value = 2421
$ws.Range("N17").Value = -2759
$ws.Range("H19").Value = 177.5
$ws.Range("I19").Value = 177.5
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 532.5
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -358.5
$ws.Range("H34").Value = 1672.5
$ws.Range("J34").Value = 4400.2
$ws.Range("L34").Value = 13200.6
$ws.Range("N34").Value = -13368.6
$ws.Range("H39").Value = 4939.4
$ws.Range("J39").Value = 5999.5
$ws.Range("L39").Value = 17998.5
$ws.Range("N39").Value = -18586.5
$ws.Range("H55").Value = 382
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("H64").Value = 4390.5
$ws.Range("J64").Value = 5123.25
$ws.Range("L64").Value = 15369.75
$ws.Range("N64").Value = -15909.75
$ws.Range("H67").Value = 4390.5
$ws.Range("J67").Value = 5123.25
$ws.Range("L67").Value = 15369.75
$ws.Range("N67").Value = -17241.75
$ws.Range("H107").Value = 877.625
$ws.Range("I107").Value = 782.2
$ws.Range("J107").Value = 1036.6666
$ws.Range("K107").Value = 2346.6
$ws.Range("L107").Value = 3109.9998
$ws.Range("M107").Value = -426.6000000000004
$ws.Range("N107").Value = -6949.9998
$ws.Range("H117").Value = 3701.5
$ws.Range("I117").Value = 714.5
$ws.Range("K117").Value = 2143.5
$ws.Range("M117").Value = 1298.5
$ws.Range("H132").Value = 37714.863
$ws.Range("I132").Value = 62587.06
$ws.Range("J132").Value = 2479.25
$ws.Range("K132").Value = 563283.54
$ws.Range("L132").Value = 22313.25
$ws.Range("M132").Value = -560753.54
$ws.Range("N132").Value = -27373.25
$ws.Range("N19").ClearContents()
$ws.Range("N55").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2342.5
$ws.Range("I102").Value = 2342.5
$ws.Range("K102").Value = 2342.5
$ws.Range("M102").Value = -720.5
$ws.Range("H126").Value = 4350.3
$ws.Range("I126").Value = 3581
$ws.Range("K126").Value = 10743
$ws.Range("M126").Value = -8273
$ws.Range("H132").Value = 6131.8335
$ws.Range("I132").Value = 4498.1055
$ws.Range("J132").Value = 12340
$ws.Range("K132").Value = 13494.3165
$ws.Range("L132").Value = 37020
$ws.Range("M132").Value = -10964.3165
$ws.Range("N132").Value = -42080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3659.1936
$ws.Range("I22").Value = 2186.5715
$ws.Range("K22").Value = 2186.5715
$ws.Range("M22").Value = -1891.5715
$ws.Range("H27").Value = 3659.1936
$ws.Range("I27").Value = 2186.5715
$ws.Range("K27").Value = 2186.5715
$ws.Range("M27").Value = -2079.5715
$ws.Range("H46").Value = 5583.727
$ws.Range("I46").Value = 1316.1111
$ws.Range("J46").Value = 8538.23
$ws.Range("K46").Value = 1316.1111
$ws.Range("L46").Value = 8538.23
$ws.Range("M46").Value = -1128.1111
$ws.Range("N46").Value = -8914.23
$ws.Range("H55").Value = 771.5
$ws.Range("I55").Value = 924.4
$ws.Range("K55").Value = 924.4
$ws.Range("M55").Value = -751.4
$ws.Range("H82").Value = 1995.6666
$ws.Range("I82").Value = 1995.6666
$ws.Range("K82").Value = 1995.6666
$ws.Range("M82").Value = -1634.6666
$ws.Range("H85").Value = 1995.6666
$ws.Range("I85").Value = 1995.6666
$ws.Range("K85").Value = 1995.6666
$ws.Range("M85").Value = -747.6666
$ws.Range("H93").Value = 2774.56
$ws.Range("I93").Value = 2540.1
$ws.Range("K93").Value = 2540.1
$ws.Range("M93").Value = -1292.1
$ws.Range("H100").Value = 6946554.5
$ws.Range("I100").Value = 9260775
$ws.Range("K100").Value = 9260775
$ws.Range("M100").Value = -9260234
$ws.Range("H132").Value = 7935.5107
$ws.Range("I132").Value = 7264.6553
$ws.Range("K132").Value = 21793.9659
$ws.Range("M132").Value = -19263.9659
$ws.Range("H136").Value = 6422.8237
$ws.Range("I136").Value = 5074.732
$ws.Range("K136").Value = 15224.196
$ws.Range("M136").Value = -12674.196

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3429.1428
$ws.Range("I132").Value = 2669.342
$ws.Range("K132").Value = 8008.026
$ws.Range("M132").Value = -5478.026
$ws.Range("H136").Value = 2282.2678
$ws.Range("I136").Value = 1442.1464
$ws.Range("J136").Value = 4578.6
$ws.Range("K136").Value = 4326.439200000001
$ws.Range("L136").Value = 13735.8
$ws.Range("M136").Value = -1776.439200000001
$ws.Range("N136").Value = -18835.8

